# Localize the "Charger sales report" workbook to Simplified Chinese:
#  - rename the worksheet tab
#  - translate the region headers that changed (Midwest/Northeast/South/Southeast)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "销售报表"

$ws.Range("B1").Value = "中西部"
$ws.Range("D1").Value = "东北"
$ws.Range("E1").Value = "南部"
$ws.Range("F1").Value = "东南部"
